$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 203, shifting existing rows 203:215 down to 204:216.
$ws.Rows.Item(203).Insert()

# Populate the newly inserted row 203 with the new record's data.
$ws.Cells.Item(203, 1).Value = 10
$ws.Cells.Item(203, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(203, 3).Value = "La Araucanía"
$ws.Cells.Item(203, 4).Value = 44610
$ws.Cells.Item(203, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(203, 5).Value = 9
$ws.Cells.Item(203, 6).Value = 100112039
$ws.Cells.Item(203, 7).Value = "Ciboulette"
$ws.Cells.Item(203, 8).Value = "Sin especificar"
$ws.Cells.Item(203, 9).Value = "Primera"
$ws.Cells.Item(203, 10).Value = 40
$ws.Cells.Item(203, 11).Value = 5000
$ws.Cells.Item(203, 12).Value = 5000
$ws.Cells.Item(203, 13).Value = 5000
$ws.Cells.Item(203, 14).Value = "`$/docena de atados"
$ws.Cells.Item(203, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(203, 16).Value = 1667
$ws.Cells.Item(203, 17).Value = 3
$ws.Cells.Item(203, 18).Value = "Hortaliza"
